# Weekly fruit/vegetable update: insert 5 new rows of Murcott mandarina
# pricing data (Región de O'Higgins, date 44491) ahead of the existing
# data block, pushing the old rows 326-365 down to 331-370.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows before row 326 - this shifts the previously
# existing rows 326:365 down to 331:370 and extends the sheet dimension
# from A1:T365 to A1:T370.
$ws.Rows("326:330").Insert()

# Constant columns shared by every data row in this sheet.
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$productoId = 100102
$producto  = "Cítricos"
$categoriaId = 100102004
$categoria = "Mandarina"
$unidad    = "$/bandeja 10 kilos"
$kgUnidad  = 10

function Set-DataRow($row, $fecha, $variedad, $calidad, $volumen, $precioMin, $precioMax, $precioProm, $origen, $precioKg) {
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

Set-DataRow 326 44491 "Murcott" "Especial" 300 5500 5500 5500 "Región de O'Higgins" 550
Set-DataRow 327 44491 "Murcott" "Extra (doble especial)" 270 6000 6000 6000 "Región de O'Higgins" 600
Set-DataRow 328 44491 "Murcott" "Primera" 380 4500 4500 4500 "Región de O'Higgins" 450
Set-DataRow 329 44491 "Murcott" "Segunda" 330 3500 3500 3500 "Región de O'Higgins" 350
Set-DataRow 330 44491 "Murcott" "Tercera" 300 2500 2500 2500 "Región de O'Higgins" 250
